# Insert a new slide (GitHub link) before the "References" slide (currently
# slide index 15), pushing References -> 16 and Thank You -> 17.
# This mirrors the sldIdLst change in the target diff:
#   ... id269, id270, id271  -->  ... id269, id272(new), id270, id271

$p = $ppt.ActivePresentation

$EMU_PER_PT = 12700

$newSlide = $p.Slides.Add(15, 12)   # 12 = ppLayoutBlank

# ---------------------------------------------------------------------
# Shape 1: "GITHUB LINK:" heading textbox
# ---------------------------------------------------------------------
$left1   = 838200   / $EMU_PER_PT
$top1    = 457200   / $EMU_PER_PT
$width1  = 6096000  / $EMU_PER_PT
$height1 = 523220   / $EMU_PER_PT

$tb1 = $newSlide.Shapes.AddTextbox(1, $left1, $top1, $width1, $height1)
$tb1.Name = "TextBox 2"
$tb1.Fill.Visible = 0

$tf1 = $tb1.TextFrame
$tf1.WordWrap = -1
$tf1.AutoSize = 1

$tr1 = $tf1.TextRange
$tr1.Text = "GITHUB LINK:"
$tr1.LanguageID = "en-US"
$tr1.Font.Size = 28
$tr1.Font.Bold = $true
$tr1.Font.Color.RGB = 8210719
$tr1.Font.Name = "+mj-lt"
$tr1.Font.NameFarEast = "Cambria"

# ---------------------------------------------------------------------
# Shape 2: hyperlinked GitHub URL textbox
# ---------------------------------------------------------------------
$left2   = 3048000  / $EMU_PER_PT
$top2    = 3246792  / $EMU_PER_PT
$width2  = 6096000  / $EMU_PER_PT
$height2 = 369332   / $EMU_PER_PT

$tb2 = $newSlide.Shapes.AddTextbox(1, $left2, $top2, $width2, $height2)
$tb2.Name = "TextBox 4"
$tb2.Fill.Visible = 0

$tf2 = $tb2.TextFrame
$tf2.WordWrap = -1
$tf2.AutoSize = 1

$tr2 = $tf2.TextRange
$tr2.Text = "https://github.com/Karthikroyal76/final-year-project"
$tr2.LanguageID = "en-IN"
$tr2.ActionSettings(1).Hyperlink.Address = "https://github.com/Karthikroyal76/final-year-project"

Write-Output "Slides.Count:"
Write-Output $p.Slides.Count
Write-Output "New slide index:"
Write-Output $newSlide.SlideIndex
Write-Output "Shapes on new slide:"
Write-Output $newSlide.Shapes.Count
